$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose data (columns B..AD) got swapped with one another.
# Column A (the running id/rank) stays fixed to its row position.
$pairs = @(
    @(18, 19),
    @(20, 21),
    @(36, 37),
    @(54, 55),
    @(58, 59),
    @(63, 64)
)

$firstCol = 2   # B
$lastCol = 30   # AD

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value()
        $v2 = $cell2.Value()

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
